$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.9583413711159094, 1.1330459575302005, 0.3648926152760323, 2.4444050437356215)"
$ws.Range("C2").Value = "NIG(1.6467694176796313, 1.2610425428217629, 3.354386461163098, 6.007340446084026)"
$ws.Range("D2").Value = "NIG(1.0174743177162338, 0.6088201534404959, 2.0246964023033565, 3.466966625511703)"
$ws.Range("E2").Value = "EXN(1.9199653738697875, 3.2805043469957873, 3.7557940208759586)"
